$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 105. Excel shifts the existing rows 105..190
# down to 106..191, carrying along their formatting (incl. the date style
# used in column D), which matches the target diff exactly.
$ws.Rows("105:105").Insert()

# Populate the newly inserted row 105 with the new data record.
$ws.Range("A105").Value2 = 7
$ws.Range("B105").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C105").Value2 = "Ñuble"
$ws.Range("D105").Value2 = 44574
$ws.Range("E105").Value2 = 16
$ws.Range("F105").Value2 = "Fruta"
$ws.Range("G105").Value2 = 100108
$ws.Range("H105").Value2 = "Tropicales y subtropicales"
$ws.Range("I105").Value2 = 100108005
$ws.Range("J105").Value2 = "Piña"
$ws.Range("K105").Value2 = "Caramelo"
$ws.Range("L105").Value2 = "Segunda"
$ws.Range("M105").Value2 = 120
$ws.Range("N105").Value2 = 16000
$ws.Range("O105").Value2 = 17000
$ws.Range("P105").Value2 = 16500
$ws.Range("Q105").Value2 = "`$/caja 14 unidades"
$ws.Range("R105").Value2 = "Ecuador"
$ws.Range("S105").Value2 = 1179
$ws.Range("T105").Value2 = 14
